# Apply cryptos list update (prices + volume % changes, Polkadot/Chainlink row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.428.43"
$ws.Range("E2").Value = "  -2.93%  "
$ws.Range("D3").Value = "3.785.99"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'598.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.88%  "
$ws.Range("D6").Value = "'167.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.72%  "
$ws.Range("D7").Value = "3.781.83"
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "'0.527"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").Value = "'0.157"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.81%  "
$ws.Range("E11").Value = "  -4.93%  "
$ws.Range("D12").Value = "'0.462"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.84%  "
$ws.Range("D13").Value = "'38.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.35%  "
$ws.Range("D14").Value = "'0.0000242"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.03%  "
$ws.Range("D15").Value = "4.425.66"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "3.799.82"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").Value = "67.535.21"
$ws.Range("E17").Value = "  -2.86%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'17.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.81%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'7.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.36%  "
$ws.Range("D21").Value = "'491.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.48%  "
$ws.Range("D22").Value = "'9.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.57%  "
$ws.Range("D23").Value = "'0.735"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Value = "'85.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("D25").Value = "'0.0000146"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.46%  "
$ws.Range("D26").Value = "'2.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.52%  "
$ws.Range("D27").Value = "'12.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.01%  "
$ws.Range("D28").Value = "'10.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.70%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -1.55%  "
$ws.Range("D31").Value = "'2.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.66%  "
$ws.Range("D32").Value = "'32.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.07%  "
$ws.Range("D33").Value = "'7.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.18%  "
$ws.Range("D34").Value = "'0.108"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.56%  "
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("D36").Value = "'1.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.41%  "
$ws.Range("E37").Value = "  -5.35%  "
$ws.Range("E38").Value = "  -4.84%  "
$ws.Range("D39").Value = "'464.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("D40").Value = "'0.326"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.99%  "
$ws.Range("D41").Value = "'49.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.13%  "
$ws.Range("D42").Value = "'1.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.61%  "
$ws.Range("E43").Value = "  -6.29%  "
$ws.Range("D44").Value = "'8.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.67%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'40.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.43%  "
$ws.Range("D47").Value = "2.832.08"
$ws.Range("E47").Value = "  -4.14%  "
$ws.Range("D48").Value = "'140.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("D49").Value = "'0.0348"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.00%  "
$ws.Range("D50").Value = "'24.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.94%  "
$ws.Range("D51").Value = "'25.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.47%  "
